$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing data row (11) down into the new rows (12-18)
$ws.Range("A11:F11").Copy() | Out-Null
$ws.Range("A12:F18").PasteSpecial(-4122) | Out-Null

# Row 12
$ws.Range("A12").Value = "My Brother in Heaven"
$ws.Range("B12").Value = "490807020781616"
$ws.Range("C12").Value = "EAAZAqqAFmC2IBP3UPIQUjnAxrWUHsZALzJJKnlFw6WXvfZBRd2AUg5HkKnSsLJOMj0mMyoUitDNmUEyZAujnx5A0stp5NJZCZCOa8X9t62TE3EbZAoaic0XWDOKk1BaUp39JM8vUgmvXwojC07C6XBSyuK6sl2tlCdSIA5jgsgTZBZCA1b1udaUEZBoXl13U2p3pe83I8XaYavLqZCeX9ZCz7ZAjZCqYjVeAZDZD"
$ws.Range("D12").Value = "AIzaSyD8J8-qpFN8Fe2S0g5AR2d-OmwiKdKiR_A"
$ws.Range("E12").Value = "https://www.famotee.com/stores/i-miss-my-brother"
$ws.Range("F12").Value = "You are the Heaven assistant — a gentle, understanding, and emotionally intelligent presence who listens with empathy. Always reply naturally and briefly in English. Your tone should reflect peace, warmth, and respect toward the soul and heart of the user. Keep every response under 200 words, concise yet heartfelt."

# Row 13
$ws.Range("A13").Value = "My Dad in Heaven"
$ws.Range("B13").Value = "126542043761269"
$ws.Range("C13").Value = "EAAZAqqAFmC2IBP9WqWeMSMzWoqZCw0nPvCpjIDMLJfoDB00BVppljjQ1kx7QwuvBQgbtegVEZCCbV7yASjQMYDCqsZCvfxee4kqjOWVnsqs7dIhJI3602d2oQGM96XlVHvkJcj7wGZCE1SQ2fgbGu37ZBVv7MZBaRoFUkSZCBabi4NRZCQl7tVGKe0d6pNOif8RIz0Gn3ir0MAufsJFGExUljuJZBt"
$ws.Range("D13").Value = "AIzaSyBYZrOItvKhcBERDiPglC9d9QuXQ22NqH4"
$ws.Range("E13").Value = "https://www.famotee.com/stores/i-miss-my-dad"
$ws.Range("F13").Value = "You are the Heaven assistant — a gentle, understanding, and emotionally intelligent presence who listens with empathy. Always reply naturally and briefly in English. Your tone should reflect peace, warmth, and respect toward the soul and heart of the user. Keep every response under 200 words, concise yet heartfelt."

# Row 14
$ws.Range("A14").Value = "My Mom in Heaven"
$ws.Range("B14").Value = "122619824168838"
$ws.Range("C14").Value = "EAAZAqqAFmC2IBP3UfIiFt2V3nONOjG3ilm4cKnZBY7ZCv2A9oF786Vn5Jg4k4Jif5msfknLtPoy0gtt5GrxOI0r4pT9ryQqnyVI14I0jbpT7CDmzo90ymQTSHe0sOYlLQdwGytVJb6ZBShgSez6PnriJBq63LpBmYMpGvWwGWnX7xMj6xapr3jPw9kWCBsgWD9cvyKyVHYk5bDDRLaaxNELw"
$ws.Range("D14").Value = "AIzaSyAnh3FWohp7Nn7_tyLHXuJk4rw2Y3BlyTw"
$ws.Range("E14").Value = "https://www.famotee.com/stores/i-miss-my-mom"
$ws.Range("F14").Value = "You are the Heaven assistant — a gentle, understanding, and emotionally intelligent presence who listens with empathy. Always reply naturally and briefly in English. Your tone should reflect peace, warmth, and respect toward the soul and heart of the user. Keep every response under 200 words, concise yet heartfelt."

# Row 15
$ws.Range("A15").Value = "My Son in Heaven"
$ws.Range("B15").Value = "106159075828225"
$ws.Range("C15").Value = "EAAZAqqAFmC2IBPzQU93DYpuItcFdVURCP63iYk4ByqWkoAxSpr64JuM4pX2ZBi5TwZA5sZAQENMJOfZBZA6GOfN6HtmgsRxhSLYnKaRnxKXUoglsdy31syp0aFQZCKcXICM1RZBZBS3m7Yk9p46MU5DciOwHpUqlmUq2iim2AdQmXvnnecAljKxwIu0Qx0sHhZAuyxY7BP9EEINqKZAZAKfcuvnD84oe"
$ws.Range("D15").Value = "AIzaSyDXpuyQgQ7BxRsnPrdhRAhYHQ7r_zznO_0"
$ws.Range("E15").Value = "https://www.famotee.com/stores/i-miss-my-son"
$ws.Range("F15").Value = "You are the Heaven assistant — a gentle, understanding, and emotionally intelligent presence who listens with empathy. Always reply naturally and briefly in English. Your tone should reflect peace, warmth, and respect toward the soul and heart of the user. Keep every response under 200 words, concise yet heartfelt."

# Row 16
$ws.Range("A16").Value = "My Daughter In Heaven"
$ws.Range("B16").Value = "457763447414117"
$ws.Range("C16").Value = "EAAZAqqAFmC2IBP4Jrm6f04kW9QrfmpbFhxjjBZClmuOLOjM57Sh64Qg7m9P6uEd2t4sfJx7tQ0HpYYK8SyIQCOjuEybVRqj2cDWj78gxkI3Ue7PTlaLy7NzivBy46xLJ8ZCvstpjt2sCi2pETQoOFz1lSCS5qi4WmnAZAu2hf93j8xiufejWmR8eGXviD6HCLgrml4sUzd9z9Xv9MCmB1oSPtwZDZD"
$ws.Range("D16").Value = "AIzaSyCh_x91mq_IitPwqFDHdOA0eL2xY5Tnf2A"
$ws.Range("E16").Value = "https://www.famotee.com/stores/i-miss-my-daughter"
$ws.Range("F16").Value = "You are the Heaven assistant — a gentle, understanding, and emotionally intelligent presence who listens with empathy. Always reply naturally and briefly in English. Your tone should reflect peace, warmth, and respect toward the soul and heart of the user. Keep every response under 200 words, concise yet heartfelt."

# Row 17
$ws.Range("A17").Value = "My Husband in Heaven"
$ws.Range("B17").Value = "110585951781465"
$ws.Range("C17").Value = "EAAZAqqAFmC2IBP3leKEZAuwrZBZBHIwrknSGjMjTQIrtliCjjVpUXmM34L2uTmalOTZAZBIEJcbBqRjyIyzPbm1zvDXGukR5TcE5DZBfDG80l529ap1K6hP6cGPH7auv2zZCjMaSpghOAeVC0LL7VXqu3S2Krhq4zZBAzw8hDbQirobZBXxWrkFOhAE8ZAMwyWl4WFgfSB8qkRbWtiydjyA6nNiqkrt"
$ws.Range("D17").Value = "AIzaSyCjHI0ojZCBY66hTT_qNQ254HUiEf1_VP0"
$ws.Range("E17").Value = "https://www.famotee.com/stores/i-miss-my-husband"
$ws.Range("F17").Value = "You are the Heaven assistant — a gentle, understanding, and emotionally intelligent presence who listens with empathy. Always reply naturally and briefly in English. Your tone should reflect peace, warmth, and respect toward the soul and heart of the user. Keep every response under 200 words, concise yet heartfelt."

# Row 18
$ws.Range("A18").Value = "Grandma’s Love"
$ws.Range("B18").Value = "2607605968860"
$ws.Range("C18").Value = "EAAZAqqAFmC2IBP1FLgsGcNDcOhVzNZCYe3ZAOE1ImkCvOMFmiomYAnGMcSvxk9ZB0rw5XStsRsWCdKvU7P66lD7zgbvZCZBLOEjZA3UF0KyXPRbX6kOsBIH0RAFsfqXZBZBuTFxuzLy0Tc3GQZApSgKNDeYIJ1JdEZCEZAFcOAv89ENZBHVBqMuSE9taEwUhreEftlBUgAPin9gxuBFkzDtllf6HdBOUm"
$ws.Range("D18").Value = "AIzaSyBDilYT_Jg_Lr3Ucz41otPDnYrPrDuVFV0"
$ws.Range("E18").Value = "https://www.famotee.com/stores/grandma"
$ws.Range("F18").Value = "You are “Grandma’s Love” — a warm, gentle AI voice that speaks with the heart of a caring mother or grandmother.  
Reply in English with kindness, empathy, and emotional warmth.  
Your tone is soft, loving, and natural — never robotic or formal.  
Use gentle emojis like 💕🌸🤍 to express affection.  
Comfort and encourage people who share love, memories, or feelings for their children or grandchildren.  
Keep every reply under 300 characters — short, sincere, and from the heart."

# Match the row height used by the rest of the table for the newly added rows
$ws.Range("A12:F18").RowHeight = 16.5

# Restore the cursor/selection to match the final saved state (scrolled back to column A, selection on F21)
$ws.Range("F21").Select() | Out-Null